$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated s_vals data (filtered save games) for rows 2-6, columns B-G.
$data = @{
    2 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182)
    3 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 1, 9.295990156953671)
    4 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182)
    5 = @(0.6753301551942219, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 0, 1.642425054193055)
    6 = @(1.459612070389937, 1.667794583268128, 337.1190423067083, 8.660232485948974, 1, 348.9066814463153)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
    $ws.Cells.Item($r, 7).Value = $vals[5]
}
